$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header/ID row) updates
$ws.Range("B1").Value2 = 15
$ws.Range("C1").Value2 = 16
$ws.Range("D1").Value2 = 15
$ws.Range("E1").Value2 = 16

# Row 2 (CON) updates - updated meanEMG legmaxROM values
$ws.Range("B2").Value2 = 108.29594798993637
$ws.Range("C2").Value2 = 105.98916852820224
$ws.Range("D2").Value2 = 107.2116753546531
$ws.Range("E2").Value2 = 107.76606483851549

# Row 3 (STR) updates - updated meanEMG legmaxROM values
$ws.Range("B3").Value2 = 107.03031794451725
$ws.Range("C3").Value2 = 104.88524901633632
$ws.Range("D3").Value2 = 107.2281202662675
$ws.Range("E3").Value2 = 108.64319819792583

# Update the active selection to reflect the newly edited range
$ws.Range("B1:E3").Select()
